$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row to append: keep the date as literal text (leading apostrophe forces
# text interpretation instead of Excel auto-converting it to a date value),
# matching the existing rows which all store plain text values.
$ws.Range("A35").Value = "'2025-09-23"
$ws.Range("B35").Value = "15:22:00"
$ws.Range("C35").Value = "1.00 EUR = 1,736.0742"
